# Weekly update: two new "Coliflor" price rows were reported for
# Vega Central Mapocho de Santiago and need to be inserted just above the
# existing row 520, pushing the rest of the table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 520-521 (existing rows 520.. shift down to 522..)
$ws.Range("520:521").Insert()

# New row 520 - "Primera" quality
$ws.Cells.Item(520, 1).Value = 9
$ws.Cells.Item(520, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(520, 3).Value = "Metropolitana"
$ws.Cells.Item(520, 4).Value = 44610
$ws.Cells.Item(520, 5).Value = 13
$ws.Cells.Item(520, 6).Value = 100112008
$ws.Cells.Item(520, 7).Value = "Coliflor"
$ws.Cells.Item(520, 8).Value = "Sin especificar"
$ws.Cells.Item(520, 9).Value = "Primera"
$ws.Cells.Item(520, 10).Value = 1600
$ws.Cells.Item(520, 11).Value = 1100
$ws.Cells.Item(520, 12).Value = 1200
$ws.Cells.Item(520, 13).Value = 1150
$ws.Cells.Item(520, 14).Value = "$/unidad"
$ws.Cells.Item(520, 15).Value = "Región Metropolitana"
$ws.Cells.Item(520, 16).Value = 1150
$ws.Cells.Item(520, 17).Value = 1
$ws.Cells.Item(520, 18).Value = "Hortaliza"

# New row 521 - "Segunda" quality
$ws.Cells.Item(521, 1).Value = 9
$ws.Cells.Item(521, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(521, 3).Value = "Metropolitana"
$ws.Cells.Item(521, 4).Value = 44610
$ws.Cells.Item(521, 5).Value = 13
$ws.Cells.Item(521, 6).Value = 100112008
$ws.Cells.Item(521, 7).Value = "Coliflor"
$ws.Cells.Item(521, 8).Value = "Sin especificar"
$ws.Cells.Item(521, 9).Value = "Segunda"
$ws.Cells.Item(521, 10).Value = 610
$ws.Cells.Item(521, 11).Value = 900
$ws.Cells.Item(521, 12).Value = 900
$ws.Cells.Item(521, 13).Value = 900
$ws.Cells.Item(521, 14).Value = "$/unidad"
$ws.Cells.Item(521, 15).Value = "Región Metropolitana"
$ws.Cells.Item(521, 16).Value = 900
$ws.Cells.Item(521, 17).Value = 1
$ws.Cells.Item(521, 18).Value = "Hortaliza"
